# Presentation state 11.02 - fix naive component forecaster bug.
#
# The naive QoQ error series is a rolling window of per-horizon errors: each
# row holds the most recent error observations for forecast horizon Q0..Q9,
# written into columns B..K (oldest->newest read left->right... actually the
# series simply grows one observation per refresh). The forecaster bug meant
# the newest observation for each horizon was never being written into the
# front of the row - this restores that: prepend the freshly computed error
# into column B of each row, shifting the previously recorded errors one
# column to the right, dropping whatever falls off the end of the window
# (column L / past column K).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Newest error observation for each row (row number -> new value for column B).
$newValues = @{
  2  = -0.5616080510579985
  3  = 0.1417647591280393
  4  = -0.4790798465348092
  5  = 0.1916007792754515
  6  = 1.573432754301089
  7  = 0.9422837133007778
  8  = 0.0678490295623069
  9  = -0.5264228954459207
  10 = 0.8949500190880419
  11 = 0.2303995154407018
  12 = 0.4008418571243615
  13 = 0.2679782848922332
  14 = -0.5417707991668423
  15 = 0.0506862842519193
  16 = -0.1624199859130616
}

# Number of previously-recorded observations (in columns B..) present in each
# row before this refresh - i.e. how far right values need to shift.
$existingCount = @{
  2  = 10
  3  = 10
  4  = 10
  5  = 10
  6  = 10
  7  = 9
  8  = 8
  9  = 7
  10 = 6
  11 = 5
  12 = 4
  13 = 3
  14 = 2
  15 = 1
  16 = 0
}

$firstCol = 2   # column B
$lastCol = 11   # column K - the window never extends past this

foreach ($row in 2..16) {
  $count = $existingCount[$row]

  # Shift existing values one column to the right, starting from the
  # rightmost populated column so we never overwrite a value before reading
  # it. Anything that would land past column K falls out of the window.
  for ($i = $count; $i -ge 1; $i--) {
    $srcCol = $firstCol + $i - 1
    $dstCol = $firstCol + $i
    $val = $ws.Cells.Item($row, $srcCol).Value2
    if ($dstCol -le $lastCol) {
      $ws.Cells.Item($row, $dstCol).Value2 = $val
    }
  }

  # Write the newly computed observation into the freed-up front column.
  $ws.Cells.Item($row, $firstCol).Value2 = $newValues[$row]
}
